# Update column F values on multiple sheets to reflect refreshed attendee/stat counts
# as generated at commit 456a3b4 (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1166
$ws.Range("F3").Value = 1079
$ws.Range("F4").Value = 1876
$ws.Range("F5").Value = 590
$ws.Range("F6").Value = 1232
$ws.Range("F7").Value = 64
$ws.Range("F8").Value = 21
$ws.Range("F10").Value = 322
$ws.Range("F11").Value = 100
$ws.Range("F12").Value = 96
$ws.Range("F13").Value = 765
$ws.Range("F14").Value = 213
$ws.Range("F15").Value = 119
$ws.Range("F18").Value = 338
$ws.Range("F19").Value = 191
$ws.Range("F20").Value = 687
$ws.Range("F21").Value = 59
$ws.Range("F23").Value = 178
$ws.Range("F24").Value = 42
$ws.Range("F25").Value = 891
$ws.Range("F26").Value = 336
$ws.Range("F27").Value = 180
$ws.Range("F28").Value = 54
$ws.Range("F32").Value = 417

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 328
$ws.Range("F6").Value = 28

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 321

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 321
$ws.Range("F3").Value = 1166
$ws.Range("F4").Value = 1079
$ws.Range("F5").Value = 1876
$ws.Range("F6").Value = 590
$ws.Range("F7").Value = 1232
$ws.Range("F8").Value = 64
$ws.Range("F10").Value = 21
$ws.Range("F12").Value = 322
$ws.Range("F13").Value = 100
$ws.Range("F14").Value = 96
$ws.Range("F15").Value = 765
$ws.Range("F16").Value = 213
$ws.Range("F17").Value = 119
$ws.Range("F20").Value = 328
$ws.Range("F23").Value = 338
$ws.Range("F24").Value = 28
$ws.Range("F27").Value = 191
$ws.Range("F28").Value = 687
$ws.Range("F29").Value = 59
$ws.Range("F31").Value = 178
$ws.Range("F32").Value = 42
$ws.Range("F33").Value = 891
$ws.Range("F34").Value = 336
$ws.Range("F37").Value = 180
$ws.Range("F38").Value = 54
$ws.Range("F46").Value = 417
